{"js": "// Apply Grammarly-style grammar/clarity fixes throughout the document.\n// Each entry is [oldText, newText]; oldText is matched literally (exact,\n// case-sensitive) against body text and is unique in the document, so a\n// body.search + Range.insertText(\"Replace\") round-trip is safe for each.\nconst changes = [\n  [\n    \"Humans are considered vertebrates meaning that we have a vertebral column, or backbone, AKA a spine.\",\n    \"Humans are considered vertebrates, meaning we have a vertebral column, or backbone, AKA a spine.\"\n  ],\n  [\n    \"AKA an articulation; The location where 2 bones come together.\",\n    \"AKA an articulation; The location where two bones come together.\"\n  ],\n  [\n    \"3 categories: immovable, slightly movable, and freely moveable\",\n    \"3 categories: immovable, slightly movable, and freely movable\"\n  ],\n  [\n    \"Tough yet flexible tissue that covers the ends of bones of the freely moveable joints.\",\n    \"Tough yet flexible tissue that covers the ends of bones of the freely movable joints.\"\n  ],\n  [\n    \"Helps Protect the bones by preventing them from rubbing together.\",\n    \"Helps protect the bones by preventing them from rubbing together.\"\n  ],\n  [\n    \"Gives support and shape to other parts of the body (e.g., ears, nose, windpipe, etc.)\",\n    \"Gives support and shape to other body parts (e.g., ears, nose, windpipe, etc.)\"\n  ],\n  [\n    \"enable the heart to beat, the chest to rise and fall, blood vessels to help regulate the pressure and flow of blood and even allow you to smile and talk.\",\n    \"enable the heart to beat, the chest to rise and fall, blood vessels to help regulate blood pressure and flow, and even allow you to smile and talk.\"\n  ],\n  [\n    \"Muscles help with movement and muscle contractions help with posture, joint stability, and heat production.\",\n    \"Muscles help with movement, and muscle contractions help with posture, joint stability, and heat production.\"\n  ],\n  [\n    \"It can be viewed as the control center for all our actions whether they are conscious or unconscious actions.\",\n    \"It can be viewed as the control center for all our actions, whether conscious or unconscious.\"\n  ],\n  [\n    \"The state of balance among all bodily systems needed for the body to survive and function correctly.\",\n    \"The state of balance among all bodily systems is needed for the body to survive and function correctly.\"\n  ],\n  [\n    \"Messages from the CNS travel through the spine to the brain which then directs the functions of the body.\",\n    \"Messages from the CNS travel through the spine to the brain, which then directs the body's functions.\"\n  ],\n  [\n    \"that can then be communicated to the brain via the cranial nerves. \",\n    \"that can be communicated to the brain via the cranial nerves. \"\n  ],\n  [\n    \"Works with the cardiovascular system to provide oxygen to cells and to remove waste products such as carbon dioxide.\",\n    \"Works with the cardiovascular system to provide cells with oxygen and remove waste products such as carbon dioxide.\"\n  ],\n  [\n    \"He primary organs that exchange gases during breathing.\",\n    \"The primary organs that exchange gases during breathing.\"\n  ],\n  [\n    \"Divided into the Upper Respiratory Tract consisting of the nose, pharynx (throat), and larynx (voice box) and the Lower Respiratory Tract consisting of the trachea (windpipe), bronchial tree, and lungs.\",\n    \"Divided into the Upper Respiratory Tract, consisting of the nose, pharynx (throat), and larynx (voice box), and the Lower Respiratory Tract, consisting of the trachea (windpipe), bronchial tree, and lungs.\"\n  ],\n  [\n    \"The process tat results in the delivery of oxygen into the body followed by the excretion of carbon dioxide outside of the body.\",\n    \"The process that results in oxygen delivery into the body, followed by the excretion of carbon dioxide outside the body.\"\n  ],\n  [\n    \"alveoli, tiny air sacks in the lungs, where after taking in oxygen the alveoli push the oxygen into the bloodstream via the capillaries. Carbon dioxide then goes through the reverse process until it is expelled via the mouth and nose.\",\n    \"alveoli, tiny air sacs in the lungs, where, after taking in oxygen, the alveoli push the oxygen into the bloodstream via the capillaries. Carbon dioxide then goes through the reverse process until expelled via the mouth and nose.\"\n  ],\n  [\n    \"body and it can be found\",\n    \"body, and it can be found\"\n  ],\n  [\n    \"The muscles on the left side of the heart pump blood through the largest single artery in the body called the aorta.\",\n    \"The muscles on the left side of the heart pump blood through the largest single artery in the body, called the aorta.\"\n  ],\n  [\n    \"Arteries get smaller as they approach the organ and at the end, the oxygen-rich blood needs to pass through the tiniest of all blood vessels called the capillary to enter the organ.\",\n    \"Arteries get smaller as they approach the organ, and at the end, the oxygen-rich blood needs to pass through the tiniest of all blood vessels, called the capillary, to enter the organ.\"\n  ],\n  [\n    \"It provides a constant supply of nutrients and oxygen to the body\u2019s cells through the flow of blood, this is known as circulation.\",\n    \"It provides a constant supply of nutrients and oxygen to the body\u2019s cells through blood flow, known as circulation.\"\n  ],\n  [\n    \"waste (carbon-dioxide and water) from the cells back to the lungs where the blood\",\n    \"waste (carbon dioxide and water) from the cells back to the lungs, where the blood\"\n  ],\n  [\n    \"Swallowed and travels down the esophagus connecting to the stomach\",\n    \"Swallowed and travels down the esophagus, connecting to the stomach\"\n  ],\n  [\n    \"small intestine where nutrient absorption occurs.\",\n    \"small intestine, where nutrient absorption occurs.\"\n  ],\n  [\n    \"exit the small intestines and go into the large intestine where water, some nutrients, and electrolytes are removed to form concentrated, solid feces.\",\n    \"exit the small intestines and go into the large intestine, where water, some nutrients, and electrolytes are removed to form concentrated, solid feces.\"\n  ],\n  [\n    \"Glands located in the mouth that moisten and lubricate food and being the breakdown of carbohydrates.\",\n    \"Glands located in the mouth that moisten and lubricate food and are involved in the breakdown of carbohydrates.\"\n  ],\n  [\n    \": Stores the yellow-green fluid called bile which helps digest and absorb fats.\",\n    \": Stores the yellow-green fluid called bile, which helps digest and absorb fats.\"\n  ],\n  [\n    \": The major bean-shaped organs of the urinary system responsible for filtering out waste products from the bloodstream and removing them as urine.\",\n    \": The major bean-shaped organs of the urinary system are responsible for filtering out waste products from the bloodstream and removing them as urine.\"\n  ],\n  [\n    \"Tubes connected to the kidneys that allow urine to flow into the urinary bladder.\",\n    \"Tubes connected to the kidneys allow urine to flow into the bladder.\"\n  ],\n  [\n    \"Protecting, Regulating, and Reproduction\",\n    \"Protecting, Regulating, and Reproducing\"\n  ]\n];\n\nfor (const [oldText, newText] of changes) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Replace only the first (and expected-only) match to stay precise even\n  // if matching text were to appear more than once.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply Grammarly-style grammar/clarity fixes throughout the document.\n# Each pair is a unique, exact sentence/phrase (oldText) located via\n# Find.Execute and then overwritten in place via Range.Text so that\n# Word's Find/Replace smart-quote AutoFormat never touches newText\n# (e.g. the literal straight apostrophe in \"body's functions\").\n$d = $word.ActiveDocument\n\n$changes = @(\n    ,@(\"Humans are considered vertebrates meaning that we have a vertebral column, or backbone, AKA a spine.\", \"Humans are considered vertebrates, meaning we have a vertebral column, or backbone, AKA a spine.\")\n    ,@(\"AKA an articulation; The location where 2 bones come together.\", \"AKA an articulation; The location where two bones come together.\")\n    ,@(\"3 categories: immovable, slightly movable, and freely moveable\", \"3 categories: immovable, slightly movable, and freely movable\")\n    ,@(\"Tough yet flexible tissue that covers the ends of bones of the freely moveable joints.\", \"Tough yet flexible tissue that covers the ends of bones of the freely movable joints.\")\n    ,@(\"Helps Protect the bones by preventing them from rubbing together.\", \"Helps protect the bones by preventing them from rubbing together.\")\n    ,@(\"Gives support and shape to other parts of the body (e.g., ears, nose, windpipe, etc.)\", \"Gives support and shape to other body parts (e.g., ears, nose, windpipe, etc.)\")\n    ,@(\"enable the heart to beat, the chest to rise and fall, blood vessels to help regulate the pressure and flow of blood and even allow you to smile and talk.\", \"enable the heart to beat, the chest to rise and fall, blood vessels to help regulate blood pressure and flow, and even allow you to smile and talk.\")\n    ,@(\"Muscles help with movement and muscle contractions help with posture, joint stability, and heat production.\", \"Muscles help with movement, and muscle contractions help with posture, joint stability, and heat production.\")\n    ,@(\"It can be viewed as the control center for all our actions whether they are conscious or unconscious actions.\", \"It can be viewed as the control center for all our actions, whether conscious or unconscious.\")\n    ,@(\"The state of balance among all bodily systems needed for the body to survive and function correctly.\", \"The state of balance among all bodily systems is needed for the body to survive and function correctly.\")\n    ,@(\"Messages from the CNS travel through the spine to the brain which then directs the functions of the body.\", \"Messages from the CNS travel through the spine to the brain, which then directs the body's functions.\")\n    ,@(\"that can then be communicated to the brain via the cranial nerves. \", \"that can be communicated to the brain via the cranial nerves. \")\n    ,@(\"Works with the cardiovascular system to provide oxygen to cells and to remove waste products such as carbon dioxide.\", \"Works with the cardiovascular system to provide cells with oxygen and remove waste products such as carbon dioxide.\")\n    ,@(\"He primary organs that exchange gases during breathing.\", \"The primary organs that exchange gases during breathing.\")\n    ,@(\"Divided into the Upper Respiratory Tract consisting of the nose, pharynx (throat), and larynx (voice box) and the Lower Respiratory Tract consisting of the trachea (windpipe), bronchial tree, and lungs.\", \"Divided into the Upper Respiratory Tract, consisting of the nose, pharynx (throat), and larynx (voice box), and the Lower Respiratory Tract, consisting of the trachea (windpipe), bronchial tree, and lungs.\")\n    ,@(\"The process tat results in the delivery of oxygen into the body followed by the excretion of carbon dioxide outside of the body.\", \"The process that results in oxygen delivery into the body, followed by the excretion of carbon dioxide outside the body.\")\n    ,@(\"alveoli, tiny air sacks in the lungs, where after taking in oxygen the alveoli push the oxygen into the bloodstream via the capillaries. Carbon dioxide then goes through the reverse process until it is expelled via the mouth and nose.\", \"alveoli, tiny air sacs in the lungs, where, after taking in oxygen, the alveoli push the oxygen into the bloodstream via the capillaries. Carbon dioxide then goes through the reverse process until expelled via the mouth and nose.\")\n    ,@(\"body and it can be found\", \"body, and it can be found\")\n    ,@(\"The muscles on the left side of the heart pump blood through the largest single artery in the body called the aorta.\", \"The muscles on the left side of the heart pump blood through the largest single artery in the body, called the aorta.\")\n    ,@(\"Arteries get smaller as they approach the organ and at the end, the oxygen-rich blood needs to pass through the tiniest of all blood vessels called the capillary to enter the organ.\", \"Arteries get smaller as they approach the organ, and at the end, the oxygen-rich blood needs to pass through the tiniest of all blood vessels, called the capillary, to enter the organ.\")\n    ,@(\"It provides a constant supply of nutrients and oxygen to the body\u2019s cells through the flow of blood, this is known as circulation.\", \"It provides a constant supply of nutrients and oxygen to the body\u2019s cells through blood flow, known as circulation.\")\n    ,@(\"waste (carbon-dioxide and water) from the cells back to the lungs where the blood\", \"waste (carbon dioxide and water) from the cells back to the lungs, where the blood\")\n    ,@(\"Swallowed and travels down the esophagus connecting to the stomach\", \"Swallowed and travels down the esophagus, connecting to the stomach\")\n    ,@(\"small intestine where nutrient absorption occurs.\", \"small intestine, where nutrient absorption occurs.\")\n    ,@(\"exit the small intestines and go into the large intestine where water, some nutrients, and electrolytes are removed to form concentrated, solid feces.\", \"exit the small intestines and go into the large intestine, where water, some nutrients, and electrolytes are removed to form concentrated, solid feces.\")\n    ,@(\"Glands located in the mouth that moisten and lubricate food and being the breakdown of carbohydrates.\", \"Glands located in the mouth that moisten and lubricate food and are involved in the breakdown of carbohydrates.\")\n    ,@(\": Stores the yellow-green fluid called bile which helps digest and absorb fats.\", \": Stores the yellow-green fluid called bile, which helps digest and absorb fats.\")\n    ,@(\": The major bean-shaped organs of the urinary system responsible for filtering out waste products from the bloodstream and removing them as urine.\", \": The major bean-shaped organs of the urinary system are responsible for filtering out waste products from the bloodstream and removing them as urine.\")\n    ,@(\"Tubes connected to the kidneys that allow urine to flow into the urinary bladder.\", \"Tubes connected to the kidneys allow urine to flow into the bladder.\")\n    ,@(\"Protecting, Regulating, and Reproduction\", \"Protecting, Regulating, and Reproducing\")\n)\n\nforeach ($pair in $changes) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false)\n    if ($found) {\n        $rng.Text = $newText\n    } else {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
